$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '25.892.69'
$ws.Range("E2").Value = '  +0.50%  '
Set-TextValue $ws.Range("D3") '1.752.12'
$ws.Range("E3").Value = '  +0.32%  '
Set-TextValue $ws.Range("D4") '1.000'
$ws.Range("E4").Value = '  +0.38%  '
Set-TextValue $ws.Range("D5") '235.97'
$ws.Range("E5").Value = '  -0.92%  '
Set-TextValue $ws.Range("D6") '0.9996'
$ws.Range("E6").Value = '  +0.32%  '
Set-TextValue $ws.Range("D7") '0.5158'
$ws.Range("E7").Value = '  +4.21%  '
Set-TextValue $ws.Range("D8") '40.38'
$ws.Range("E8").Value = '  -2.97%  '
Set-TextValue $ws.Range("D9") '0.2690'
$ws.Range("E9").Value = '  +1.35%  '
Set-TextValue $ws.Range("D10") '0.06183'
$ws.Range("E10").Value = '  +1.24%  '
Set-TextValue $ws.Range("D11") '1.755.70'
$ws.Range("E11").Value = '  +0.62%  '
Set-TextValue $ws.Range("D12") '0.06977'
$ws.Range("E12").Value = '  +1.13%  '
Set-TextValue $ws.Range("D13") '15.40'
$ws.Range("E13").Value = '  +0.53%  '
Set-TextValue $ws.Range("D14") '0.6351'
$ws.Range("E14").Value = '  +8.63%  '
Set-TextValue $ws.Range("D15") '4.478'
$ws.Range("E15").Value = '  +0.05%  '
Set-TextValue $ws.Range("D16") '77.85'
$ws.Range("E16").Value = '  +1.38%  '
Set-TextValue $ws.Range("D17") '1.0000'
$ws.Range("E17").Value = '  +0.45%  '
Set-TextValue $ws.Range("D18") '0.9988'
$ws.Range("E18").Value = '  +0.07%  '
Set-TextValue $ws.Range("D19") '25.916.16'
$ws.Range("E19").Value = '  +0.47%  '
Set-TextValue $ws.Range("D20") '11.62'
$ws.Range("E20").Value = '  +0.29%  '
Set-TextValue $ws.Range("D21") '0.000006670'
$ws.Range("E21").Value = '  -0.17%  '
Set-TextValue $ws.Range("D22") '1.980.43'
$ws.Range("E22").Value = '  +0.95%  '
Set-TextValue $ws.Range("D23") '4.064'
$ws.Range("E23").Value = '  +0.42%  '
Set-TextValue $ws.Range("D24") '8.307'
$ws.Range("E24").Value = '  +4.06%  '
Set-TextValue $ws.Range("D25") '5.173'
$ws.Range("E25").Value = '  +2.04%  '
Set-TextValue $ws.Range("D26") '136.19'
$ws.Range("E26").Value = '  -1.02%  '
Set-TextValue $ws.Range("D27") '1.484'
$ws.Range("E27").Value = '  -2.50%  '
Set-TextValue $ws.Range("D28") '15.12'
$ws.Range("E28").Value = '  +1.94%  '
Set-TextValue $ws.Range("D29") '1.802'
$ws.Range("E29").Value = '  -1.98%  '
Set-TextValue $ws.Range("D30") '102.97'
$ws.Range("E30").Value = '  +0.87%  '
Set-TextValue $ws.Range("D31") '0.08303'
$ws.Range("E31").Value = '  +3.51%  '
Set-TextValue $ws.Range("D32") '3.690'
$ws.Range("E32").Value = '  -1.46%  '
Set-TextValue $ws.Range("D33") '3.394'
$ws.Range("E33").Value = '  -2.27%  '
Set-TextValue $ws.Range("D34") '0.04383'
$ws.Range("E34").Value = '  -1.34%  '
Set-TextValue $ws.Range("D35") '2.639'
$ws.Range("E35").Value = '  +0.61%  '
Set-TextValue $ws.Range("D36") '0.9942'
$ws.Range("E36").Value = '  +1.53%  '
Set-TextValue $ws.Range("D37") '0.6008'
$ws.Range("E37").Value = '  +0.33%  '
Set-TextValue $ws.Range("D39") '0.01559'
$ws.Range("E39").Value = '  +2.78%  '
Set-TextValue $ws.Range("D40") '1.925'
$ws.Range("E40").Value = '  -0.10%  '
Set-TextValue $ws.Range("D41") '0.9994'
$ws.Range("E41").Value = '  +0.28%  '
Set-TextValue $ws.Range("D42") '102.15'
$ws.Range("E42").Value = '  -2.71%  '
Set-TextValue $ws.Range("D43") '0.3852'
$ws.Range("E43").Value = '  +1.53%  '
Set-TextValue $ws.Range("D44") '0.7479'
$ws.Range("E44").Value = '  +3.12%  '
Set-TextValue $ws.Range("D45") '4.905'
$ws.Range("E45").Value = '  -4.75%  '
Set-TextValue $ws.Range("D47") '0.1103'
$ws.Range("E47").Value = '  -0.57%  '
Set-TextValue $ws.Range("D48") '6.010'
$ws.Range("E48").Value = '  +1.79%  '
Set-TextValue $ws.Range("D49") '30.21'
$ws.Range("E49").Value = '  +0.44%  '
Set-TextValue $ws.Range("D50") '52.53'
$ws.Range("E50").Value = '  +0.57%  '
Set-TextValue $ws.Range("D51") '1.004'
$ws.Range("E51").Value = '  +0.73%  '
$ws.Range("E38").Value = '  +2.49%  '
$ws.Range("E46").Value = '  +5.71%  '
